# chore: update Sheets via scheduled runner
# Refresh of market-board derived profit columns (H,I,J,K,L,M,N) for select
# leve rows across the crafting job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5084.5
$ws.Range("I32").Value = 2532.3333
$ws.Range("K32").Value = 2532.3333
$ws.Range("M32").Value = -2206.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2632.1738
$ws.Range("J51").Value = 2934.2856
$ws.Range("L51").Value = 2934.2856
$ws.Range("N51").Value = -3902.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5980.6
$ws.Range("I98").Value = 5999.3335
$ws.Range("J98").Value = 5952.5
$ws.Range("K98").Value = 5999.3335
$ws.Range("L98").Value = 5952.5
$ws.Range("M98").Value = -4501.3335
$ws.Range("N98").Value = -8948.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 5980.6
$ws.Range("I122").Value = 5999.3335
$ws.Range("J122").Value = 5952.5
$ws.Range("K122").Value = 17998.0005
$ws.Range("L122").Value = 17857.5
$ws.Range("M122").Value = -15548.0005
$ws.Range("N122").Value = -22757.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1886.9688
$ws.Range("I137").Value = 1823.4706
$ws.Range("J137").Value = 1958.9333
$ws.Range("K137").Value = 5470.4118
$ws.Range("L137").Value = 5876.7999
$ws.Range("M137").Value = -2920.4118
$ws.Range("N137").Value = -10976.7999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5497192.5
$ws.Range("I138").Value = 1229.76
$ws.Range("J138").Value = 7578996.5
$ws.Range("K138").Value = 3689.28
$ws.Range("L138").Value = 22736989.5
$ws.Range("M138").Value = 1450.72
$ws.Range("N138").Value = -22747269.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7684.9
$ws.Range("I32").Value = 4339.6416
$ws.Range("K32").Value = 4339.6416
$ws.Range("M32").Value = -4052.6416

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3216.1177
$ws.Range("I61").Value = 2619.5715
$ws.Range("K61").Value = 2619.5715
$ws.Range("M61").Value = -2407.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6710.3887
$ws.Range("I74").Value = 993.5
$ws.Range("J74").Value = 35294.832
$ws.Range("K74").Value = 993.5
$ws.Range("L74").Value = 35294.832
$ws.Range("M74").Value = -119.5
$ws.Range("N74").Value = -37042.832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6710.3887
$ws.Range("I77").Value = 993.5
$ws.Range("J77").Value = 35294.832
$ws.Range("K77").Value = 4967.5
$ws.Range("L77").Value = 176474.16
$ws.Range("M77").Value = -599.5
$ws.Range("N77").Value = -185210.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3216.1177
$ws.Range("I136").Value = 2619.5715
$ws.Range("K136").Value = 7858.7145
$ws.Range("M136").Value = -5308.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3314.7
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1733.5
$ws.Range("I134").Value = 1733.5
$ws.Range("K134").Value = 5200.5
$ws.Range("M134").Value = -2665.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 103108
$ws.Range("I31").Value = 145044.42
$ws.Range("J31").Value = 5256.3335
$ws.Range("K31").Value = 145044.42
$ws.Range("L31").Value = 5256.3335
$ws.Range("M31").Value = -144749.42
$ws.Range("N31").Value = -5846.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 103108
$ws.Range("I34").Value = 145044.42
$ws.Range("J34").Value = 5256.3335
$ws.Range("K34").Value = 145044.42
$ws.Range("L34").Value = 5256.3335
$ws.Range("M34").Value = -144842.42
$ws.Range("N34").Value = -5660.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2685.3572
$ws.Range("I99").Value = 2554.0908
$ws.Range("J99").Value = 3166.6667
$ws.Range("K99").Value = 2554.0908
$ws.Range("L99").Value = 3166.6667
$ws.Range("M99").Value = -1056.0908
$ws.Range("N99").Value = -6162.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2685.3572
$ws.Range("I126").Value = 2554.0908
$ws.Range("J126").Value = 3166.6667
$ws.Range("K126").Value = 7662.2724
$ws.Range("L126").Value = 9500.000100000001
$ws.Range("M126").Value = -5192.2724
$ws.Range("N126").Value = -14440.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2612.2593
$ws.Range("I132").Value = 2443.4583
$ws.Range("K132").Value = 7330.374899999999
$ws.Range("M132").Value = -4800.374899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6231.4062
$ws.Range("I134").Value = 3532.0688
$ws.Range("K134").Value = 10596.2064
$ws.Range("M134").Value = -8061.206399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 498064.6
$ws.Range("J141").Value = 498064.6
$ws.Range("L141").Value = 498064.6
$ws.Range("N141").Value = -508424.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 676.9091
$ws.Range("I5").Value = 637
$ws.Range("K5").Value = 1911
$ws.Range("M5").Value = -1799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 14178675
$ws.Range("I9").Value = 49500036
$ws.Range("J9").Value = 50131
$ws.Range("K9").Value = 148500108
$ws.Range("L9").Value = 150393
$ws.Range("M9").Value = -148499884
$ws.Range("N9").Value = -150841

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 65.666664
$ws.Range("I10").Value = 65.666664
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 196.999992
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -57.99999199999999
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 142
$ws.Range("I33").Value = 38.25
$ws.Range("J33").Value = 319.85715
$ws.Range("K33").Value = 229.5
$ws.Range("L33").Value = 1919.1429
$ws.Range("M33").Value = 53.5
$ws.Range("N33").Value = -2485.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3224.75
$ws.Range("I116").Value = 966.3333
$ws.Range("K116").Value = 2898.9999
$ws.Range("M116").Value = 543.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22353.572
$ws.Range("I131").Value = 250475
$ws.Range("J131").Value = 2076.111
$ws.Range("K131").Value = 751425
$ws.Range("L131").Value = 6228.333
$ws.Range("M131").Value = -746385
$ws.Range("N131").Value = -16308.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 676.9091
$ws.Range("I135").Value = 637
$ws.Range("K135").Value = 5733
$ws.Range("M135").Value = -3198

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6916.5
$ws.Range("I70").Value = 6874.75
$ws.Range("K70").Value = 6874.75
$ws.Range("M70").Value = -6604.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6916.5
$ws.Range("I73").Value = 6874.75
$ws.Range("K73").Value = 6874.75
$ws.Range("M73").Value = -5938.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8126.5713
$ws.Range("I7").Value = 9834.666999999999
$ws.Range("K7").Value = 9834.666999999999
$ws.Range("M7").Value = -9722.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3129.3333
$ws.Range("I40").Value = 2738.0476
$ws.Range("K40").Value = 2738.0476
$ws.Range("M40").Value = -2602.0476

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 39812
$ws.Range("J108").Value = 39812
$ws.Range("L108").Value = 39812
$ws.Range("N108").Value = -47492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6148.154
$ws.Range("I122").Value = 5436.8887
$ws.Range("K122").Value = 16310.6661
$ws.Range("M122").Value = -13860.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8126.5713
$ws.Range("I126").Value = 9834.666999999999
$ws.Range("K126").Value = 29504.001
$ws.Range("M126").Value = -27034.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3569.0417
$ws.Range("I132").Value = 2882.95
$ws.Range("K132").Value = 8648.849999999999
$ws.Range("M132").Value = -6118.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1555.3077
$ws.Range("I122").Value = 1523.3914
$ws.Range("K122").Value = 4570.174199999999
$ws.Range("M122").Value = -2120.174199999999
